$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44294
$ws.Range("K2").Value = 14000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 14500
$ws.Range("O2").Value = "Región de Arica y Parinacota"
$ws.Range("P2").Value = 242

$ws.Range("D3").Value = 44511
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 8000
$ws.Range("L3").Value = 9000
$ws.Range("M3").Value = 8500
$ws.Range("O3").Value = "Región de Arica y Parinacota"
$ws.Range("P3").Value = 142

$ws.Range("D4").Value = 44572
$ws.Range("J4").Value = 310
$ws.Range("K4").Value = 5500
$ws.Range("L4").Value = 6000
$ws.Range("M4").Value = 5742
$ws.Range("P4").Value = 96

$ws.Range("D5").Value = 44223
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 9000
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = 9500
$ws.Range("P5").Value = 158

$ws.Range("D6").Value = 44210
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 11000
$ws.Range("L6").Value = 12000
$ws.Range("M6").Value = 11500
$ws.Range("N6").Value = "`$/caja 60 unidades"
$ws.Range("O6").Value = "Región de Arica y Parinacota"
$ws.Range("P6").Value = 192
$ws.Range("Q6").Value = 60

$ws.Range("D7").Value = 44435
$ws.Range("K7").Value = 14000
$ws.Range("L7").Value = 15000
$ws.Range("M7").Value = 14500
$ws.Range("N7").Value = "`$/caja 60 unidades"
$ws.Range("P7").Value = 242
$ws.Range("Q7").Value = 60

$ws.Range("D8").Value = 44320
$ws.Range("K8").Value = 9000
$ws.Range("L8").Value = 10000
$ws.Range("M8").Value = 9500
$ws.Range("P8").Value = 158

$ws.Range("D9").Value = 44336
$ws.Range("K9").Value = 10000
$ws.Range("L9").Value = 11000
$ws.Range("M9").Value = 10500
$ws.Range("P9").Value = 175

$ws.Range("D10").Value = 44425

$ws.Range("D11").Value = 44537
$ws.Range("J11").Value = 220
$ws.Range("K11").Value = 9000
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = 9545
$ws.Range("P11").Value = 159

$ws.Range("D12").Value = 44496
$ws.Range("J12").Value = 350
$ws.Range("M12").Value = 6786
$ws.Range("O12").Value = "Región de Coquimbo"

$ws.Range("D13").Value = 44299
$ws.Range("K13").Value = 14000
$ws.Range("L13").Value = 15000
$ws.Range("M13").Value = 14500
$ws.Range("P13").Value = 242

$ws.Range("D14").Value = 44509
$ws.Range("K14").Value = 6000
$ws.Range("L14").Value = 6500
$ws.Range("M14").Value = 6250
$ws.Range("P14").Value = 104

$ws.Range("D15").Value = 44265
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 13000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = 14000
$ws.Range("P15").Value = 233

$ws.Range("D16").Value = 44638
$ws.Range("J16").Value = 270
$ws.Range("K16").Value = 18000
$ws.Range("L16").Value = 19000
$ws.Range("M16").Value = 18556
$ws.Range("P16").Value = 309

$ws.Range("D17").Value = 44397
$ws.Range("K17").Value = 17000
$ws.Range("L17").Value = 18000
$ws.Range("M17").Value = 17500
$ws.Range("P17").Value = 292

$ws.Range("D18").Value = 44168
$ws.Range("K18").Value = 6500
$ws.Range("L18").Value = 7000
$ws.Range("M18").Value = 6750
$ws.Range("P18").Value = 112

$ws.Range("D19").Value = 44678
$ws.Range("K19").Value = 17000
$ws.Range("L19").Value = 18000
$ws.Range("M19").Value = 17500
$ws.Range("P19").Value = 292

$ws.Range("D20").Value = 44460
$ws.Range("K20").Value = 16000
$ws.Range("L20").Value = 17000
$ws.Range("M20").Value = 16500
$ws.Range("P20").Value = 275

$ws.Range("D21").Value = 44355
$ws.Range("K21").Value = 10000
$ws.Range("L21").Value = 11000
$ws.Range("M21").Value = 10500
$ws.Range("P21").Value = 175

$ws.Range("D22").Value = 44623
$ws.Range("J22").Value = 220
$ws.Range("M22").Value = 14455
$ws.Range("P22").Value = 241

$ws.Range("D23").Value = 44482
$ws.Range("J23").Value = 350
$ws.Range("K23").Value = 10000
$ws.Range("L23").Value = 11000
$ws.Range("M23").Value = 10429
$ws.Range("P23").Value = 174

$ws.Range("D24").Value = 44292
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("J24").Value = 100
$ws.Range("K24").Value = 14000
$ws.Range("L24").Value = 15000
$ws.Range("M24").Value = 14500
$ws.Range("P24").Value = 242

$ws.Range("D25").Value = 44313
$ws.Range("K25").Value = 9000
$ws.Range("L25").Value = 10000
$ws.Range("M25").Value = 9500
$ws.Range("P25").Value = 158

$ws.Range("D26").Value = 44253
$ws.Range("J26").Value = 200
$ws.Range("K26").Value = 9000
$ws.Range("L26").Value = 10000
$ws.Range("M26").Value = 9500
$ws.Range("P26").Value = 158

$ws.Range("D27").Value = 44194
$ws.Range("J27").Value = 100
$ws.Range("M27").Value = 11500
$ws.Range("P27").Value = 192

$ws.Range("D28").Value = 44285
$ws.Range("K28").Value = 12000
$ws.Range("L28").Value = 13000
$ws.Range("M28").Value = 12500
$ws.Range("P28").Value = 208

$ws.Range("D29").Value = 44392
$ws.Range("J29").Value = 100
$ws.Range("K29").Value = 16000
$ws.Range("L29").Value = 17000
$ws.Range("M29").Value = 16500
$ws.Range("P29").Value = 275

$ws.Range("D30").Value = 44258
$ws.Range("J30").Value = 200
$ws.Range("K30").Value = 12000
$ws.Range("L30").Value = 13000
$ws.Range("M30").Value = 12500
$ws.Range("P30").Value = 208

$ws.Range("D31").Value = 44203
$ws.Range("J31").Value = 100
$ws.Range("M31").Value = 11500
$ws.Range("P31").Value = 192

$ws.Range("D32").Value = 44680
$ws.Range("L32").Value = 16000
$ws.Range("M32").Value = 15500
$ws.Range("P32").Value = 258

$ws.Range("D33").Value = 44455
$ws.Range("K33").Value = 15000
$ws.Range("L33").Value = 16000
$ws.Range("M33").Value = 15500
$ws.Range("N33").Value = "`$/caja 50 unidades"
$ws.Range("P33").Value = 310
$ws.Range("Q33").Value = 50

$ws.Range("D34").Value = 44670
$ws.Range("K34").Value = 15000
$ws.Range("L34").Value = 16000
$ws.Range("M34").Value = 15455
$ws.Range("O34").Value = "Limache"
$ws.Range("P34").Value = 258

$ws.Range("D35").Value = 44484
$ws.Range("J35").Value = 450
$ws.Range("K35").Value = 11000
$ws.Range("L35").Value = 12000
$ws.Range("M35").Value = 11556
$ws.Range("O35").Value = "Región de Arica y Parinacota"
$ws.Range("P35").Value = 193

$ws.Range("D36").Value = 44624
$ws.Range("J36").Value = 270
$ws.Range("K36").Value = 14000
$ws.Range("L36").Value = 15000
$ws.Range("M36").Value = 14444
$ws.Range("O36").Value = "Región Metropolitana"
$ws.Range("P36").Value = 241

$ws.Range("D37").Value = 44608
$ws.Range("J37").Value = 100
$ws.Range("M37").Value = 14500
$ws.Range("P37").Value = 242

$ws.Range("D38").Value = 44211
$ws.Range("J38").Value = 200
$ws.Range("K38").Value = 11000
$ws.Range("L38").Value = 12000
$ws.Range("M38").Value = 11500
$ws.Range("P38").Value = 192

$ws.Range("D39").Value = 44379
$ws.Range("K39").Value = 13000
$ws.Range("L39").Value = 14000
$ws.Range("M39").Value = 13500
$ws.Range("P39").Value = 225

$ws.Range("D40").Value = 44278
$ws.Range("J40").Value = 100
$ws.Range("K40").Value = 11000
$ws.Range("L40").Value = 12000
$ws.Range("M40").Value = 11500
$ws.Range("O40").Value = "Región de Arica y Parinacota"
$ws.Range("P40").Value = 192

$ws.Range("D41").Value = 44663
$ws.Range("K41").Value = 14000
$ws.Range("L41").Value = 15000
$ws.Range("M41").Value = 14500
$ws.Range("P41").Value = 242

$ws.Range("D42").Value = 44659
$ws.Range("H42").Value = "Alaska"
$ws.Range("J42").Value = 250
$ws.Range("L42").Value = 12000
$ws.Range("M42").Value = 10800
$ws.Range("N42").Value = "`$/caja 60 unidades"
$ws.Range("O42").Value = "Región de Arica y Parinacota"
$ws.Range("P42").Value = 180
$ws.Range("Q42").Value = 60

$ws.Range("D43").Value = 44467
$ws.Range("J43").Value = 100
$ws.Range("M43").Value = 15500
$ws.Range("O43").Value = "Región de Arica y Parinacota"

$ws.Range("D44").Value = 44385
$ws.Range("J44").Value = 100
$ws.Range("K44").Value = 15000
$ws.Range("L44").Value = 16000
$ws.Range("M44").Value = 15500
$ws.Range("P44").Value = 258

$ws.Range("D45").Value = 44532
$ws.Range("K45").Value = 6500
$ws.Range("L45").Value = 7000
$ws.Range("M45").Value = 6700
$ws.Range("P45").Value = 112

$ws.Range("D46").Value = 44306
$ws.Range("J46").Value = 200
$ws.Range("K46").Value = 9000
$ws.Range("L46").Value = 10000
$ws.Range("M46").Value = 9500
$ws.Range("P46").Value = 158

$ws.Range("D47").Value = 44580
$ws.Range("J47").Value = 150
$ws.Range("K47").Value = 11000
$ws.Range("L47").Value = 12000
$ws.Range("M47").Value = 11667
$ws.Range("P47").Value = 194

$ws.Range("D48").Value = 44420
$ws.Range("J48").Value = 200
$ws.Range("K48").Value = 16000
$ws.Range("L48").Value = 17000
$ws.Range("M48").Value = 16500
$ws.Range("P48").Value = 275

$ws.Range("D49").Value = 44272
$ws.Range("J49").Value = 100
$ws.Range("K49").Value = 12000
$ws.Range("L49").Value = 13000
$ws.Range("M49").Value = 12500
$ws.Range("P49").Value = 208

$ws.Range("D50").Value = 44645
$ws.Range("J50").Value = 170
$ws.Range("K50").Value = 17000
$ws.Range("L50").Value = 18000
$ws.Range("M50").Value = 17529
$ws.Range("P50").Value = 292

$ws.Range("D51").Value = 44308
$ws.Range("K51").Value = 11000
$ws.Range("L51").Value = 12000
$ws.Range("M51").Value = 11500
$ws.Range("P51").Value = 192

$ws.Range("D52").Value = 44421
$ws.Range("J52").Value = 100
$ws.Range("K52").Value = 17000
$ws.Range("L52").Value = 18000
$ws.Range("M52").Value = 17500
$ws.Range("P52").Value = 292

$ws.Range("D53").Value = 44432
$ws.Range("K53").Value = 14000
$ws.Range("L53").Value = 15000
$ws.Range("M53").Value = 14500
$ws.Range("P53").Value = 242

$ws.Range("D54").Value = 44295
$ws.Range("K54").Value = 13000
$ws.Range("L54").Value = 14000
$ws.Range("M54").Value = 13500
$ws.Range("P54").Value = 225

$ws.Range("D55").Value = 44642
$ws.Range("J55").Value = 270
$ws.Range("K55").Value = 17000
$ws.Range("L55").Value = 18000
$ws.Range("M55").Value = 17556
$ws.Range("P55").Value = 293

$ws.Range("D56").Value = 44230
$ws.Range("K56").Value = 9000
$ws.Range("L56").Value = 10000
$ws.Range("M56").Value = 9500
$ws.Range("P56").Value = 158

$ws.Range("D57").Value = 44476
$ws.Range("K57").Value = 16000
$ws.Range("L57").Value = 17000
$ws.Range("M57").Value = 16500
$ws.Range("P57").Value = 275

$ws.Range("D58").Value = 44526
$ws.Range("J58").Value = 200

$ws.Range("D59").Value = 44350
$ws.Range("J59").Value = 100
$ws.Range("K59").Value = 10000
$ws.Range("L59").Value = 12000
$ws.Range("M59").Value = 11000
$ws.Range("P59").Value = 183

$ws.Range("D60").Value = 44530
$ws.Range("J60").Value = 350
$ws.Range("K60").Value = 6000
$ws.Range("L60").Value = 6500
$ws.Range("M60").Value = 6286
$ws.Range("N60").Value = "`$/caja 80 unidades"
$ws.Range("O60").Value = "Región del Maule"
$ws.Range("P60").Value = 79
$ws.Range("Q60").Value = 80

$ws.Range("D61").Value = 44567
$ws.Range("K61").Value = 7000
$ws.Range("L61").Value = 7500
$ws.Range("M61").Value = 7250
$ws.Range("P61").Value = 121

$ws.Range("D62").Value = 44370
$ws.Range("K62").Value = 15000
$ws.Range("L62").Value = 16000
$ws.Range("M62").Value = 15500
$ws.Range("O62").Value = "Región Metropolitana"
$ws.Range("P62").Value = 258

$ws.Range("D63").Value = 44327
$ws.Range("J63").Value = 100

$ws.Range("D64").Value = 44246
$ws.Range("J64").Value = 200
$ws.Range("K64").Value = 10000
$ws.Range("L64").Value = 12000
$ws.Range("M64").Value = 11000
$ws.Range("O64").Value = "Región del Maule"
$ws.Range("P64").Value = 183

$ws.Range("D65").Value = 44383
$ws.Range("K65").Value = 14000
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = 14500
$ws.Range("P65").Value = 242

$ws.Range("D66").Value = 44362
$ws.Range("K66").Value = 12000
$ws.Range("L66").Value = 13000
$ws.Range("M66").Value = 12500
$ws.Range("N66").Value = "`$/caja 50 unidades"
$ws.Range("P66").Value = 250
$ws.Range("Q66").Value = 50

$ws.Range("D67").Value = 44517
$ws.Range("J67").Value = 250
$ws.Range("K67").Value = 5500
$ws.Range("L67").Value = 6000
$ws.Range("M67").Value = 5700
$ws.Range("P67").Value = 95

$ws.Range("D68").Value = 44390
$ws.Range("K68").Value = 16000
$ws.Range("L68").Value = 17000
$ws.Range("M68").Value = 16500
$ws.Range("P68").Value = 275

$ws.Range("D69").Value = 44237
$ws.Range("K69").Value = 10000
$ws.Range("L69").Value = 11000
$ws.Range("M69").Value = 10500
$ws.Range("P69").Value = 175

$ws.Range("D70").Value = 44644
$ws.Range("J70").Value = 270
$ws.Range("K70").Value = 17000
$ws.Range("L70").Value = 18000
$ws.Range("M70").Value = 17444
$ws.Range("P70").Value = 291

$ws.Range("D71").Value = 44539
$ws.Range("J71").Value = 180
$ws.Range("K71").Value = 6500
$ws.Range("L71").Value = 7000
$ws.Range("M71").Value = 6722
$ws.Range("P71").Value = 112

$ws.Range("D72").Value = 44162
$ws.Range("J72").Value = 200
$ws.Range("K72").Value = 7000
$ws.Range("L72").Value = 7500
$ws.Range("M72").Value = 7250
$ws.Range("O72").Value = "Región del Maule"
$ws.Range("P72").Value = 121

$ws.Range("D73").Value = 44643
$ws.Range("J73").Value = 140
$ws.Range("K73").Value = 16000
$ws.Range("L73").Value = 17000
$ws.Range("M73").Value = 16571
$ws.Range("N73").Value = "`$/caja 70 unidades"
$ws.Range("O73").Value = "Provincia de Limarí"
$ws.Range("P73").Value = 237
$ws.Range("Q73").Value = 70

$ws.Range("D74").Value = 44196
$ws.Range("K74").Value = 11000
$ws.Range("L74").Value = 12000
$ws.Range("M74").Value = 11500
$ws.Range("P74").Value = 192

$ws.Range("D75").Value = 44453
$ws.Range("J75").Value = 100
$ws.Range("K75").Value = 16000
$ws.Range("L75").Value = 17000
$ws.Range("M75").Value = 16500
$ws.Range("P75").Value = 275

$ws.Range("D76").Value = 44503
$ws.Range("J76").Value = 250
$ws.Range("K76").Value = 7500
$ws.Range("L76").Value = 8000
$ws.Range("M76").Value = 7700
$ws.Range("P76").Value = 128

$ws.Range("D77").Value = 44628
$ws.Range("J77").Value = 220
$ws.Range("K77").Value = 17000
$ws.Range("L77").Value = 19000
$ws.Range("M77").Value = 17909
$ws.Range("P77").Value = 298

$ws.Range("D78").Value = 44341
$ws.Range("J78").Value = 100
$ws.Range("K78").Value = 9000
$ws.Range("L78").Value = 10000
$ws.Range("M78").Value = 9500
$ws.Range("P78").Value = 158

$ws.Range("D79").Value = 44491
$ws.Range("K79").Value = 8500
$ws.Range("L79").Value = 9000
$ws.Range("M79").Value = 8750
$ws.Range("P79").Value = 146

$ws.Range("D80").Value = 44685
$ws.Range("J80").Value = 220
$ws.Range("K80").Value = 17000
$ws.Range("L80").Value = 18000
$ws.Range("M80").Value = 17455
$ws.Range("O80").Value = "Región Metropolitana"
$ws.Range("P80").Value = 291

$ws.Range("D81").Value = 44330
$ws.Range("K81").Value = 10000
$ws.Range("L81").Value = 11000
$ws.Range("M81").Value = 10500
$ws.Range("P81").Value = 175

$ws.Range("D82").Value = 44399
$ws.Range("J82").Value = 100
$ws.Range("K82").Value = 16000
$ws.Range("L82").Value = 17000
$ws.Range("M82").Value = 16500
$ws.Range("N82").Value = "`$/caja 60 unidades"
$ws.Range("O82").Value = "Región de Arica y Parinacota"
$ws.Range("P82").Value = 275
$ws.Range("Q82").Value = 60

$ws.Range("D83").Value = 44238
$ws.Range("K83").Value = 12000
$ws.Range("L83").Value = 14000
$ws.Range("M83").Value = 13000
$ws.Range("O83").Value = "Región de Arica y Parinacota"
$ws.Range("P83").Value = 217

$ws.Range("D84").Value = 44250
$ws.Range("J84").Value = 100
$ws.Range("K84").Value = 8000
$ws.Range("L84").Value = 9000
$ws.Range("M84").Value = 8500
$ws.Range("P84").Value = 142

$ws.Range("D85").Value = 44334
$ws.Range("K85").Value = 11000
$ws.Range("L85").Value = 12000
$ws.Range("M85").Value = 11500
$ws.Range("P85").Value = 192

$ws.Range("D86").Value = 44622
$ws.Range("J86").Value = 250
$ws.Range("K86").Value = 8500
$ws.Range("L86").Value = 9000
$ws.Range("M86").Value = 8800
$ws.Range("P86").Value = 147

$ws.Range("D87").Value = 44637
$ws.Range("J87").Value = 270
$ws.Range("K87").Value = 17000
$ws.Range("L87").Value = 18000
$ws.Range("M87").Value = 17556
$ws.Range("P87").Value = 293

$ws.Range("D88").Value = 44498
$ws.Range("K88").Value = 6500
$ws.Range("L88").Value = 7000
$ws.Range("M88").Value = 6786
$ws.Range("P88").Value = 113

$ws.Range("D89").Value = 44189
$ws.Range("K89").Value = 11000
$ws.Range("L89").Value = 12000
$ws.Range("M89").Value = 11500
$ws.Range("P89").Value = 192

$ws.Range("D90").Value = 44217
$ws.Range("J90").Value = 200
$ws.Range("K90").Value = 8000
$ws.Range("L90").Value = 9000
$ws.Range("M90").Value = 8500
$ws.Range("N90").Value = "`$/caja 60 unidades"
$ws.Range("O90").Value = "Región del Maule"
$ws.Range("P90").Value = 142
$ws.Range("Q90").Value = 60

$ws.Range("D91").Value = 44215
$ws.Range("J91").Value = 100
$ws.Range("K91").Value = 9000
$ws.Range("L91").Value = 10000
$ws.Range("M91").Value = 9500
$ws.Range("P91").Value = 158

$ws.Range("D92").Value = 44358

$ws.Range("D93").Value = 44658
$ws.Range("J93").Value = 220
$ws.Range("K93").Value = 14000
$ws.Range("L93").Value = 15000
$ws.Range("M93").Value = 14545
$ws.Range("N93").Value = "`$/caja 70 unidades"
$ws.Range("O93").Value = "Región del Maule"
$ws.Range("P93").Value = 208
$ws.Range("Q93").Value = 70

$ws.Range("D94").Value = 44551
$ws.Range("K94").Value = 7000
$ws.Range("L94").Value = 8000
$ws.Range("M94").Value = 7500
$ws.Range("P94").Value = 125

$ws.Range("D95").Value = 44281
$ws.Range("K95").Value = 12000
$ws.Range("L95").Value = 13000
$ws.Range("M95").Value = 12500
$ws.Range("P95").Value = 208

$ws.Range("D96").Value = 44691
$ws.Range("K96").Value = 16000
$ws.Range("L96").Value = 17000
$ws.Range("M96").Value = 16500
$ws.Range("P96").Value = 275

$ws.Range("D97").Value = 44187
$ws.Range("J97").Value = 200
$ws.Range("K97").Value = 8000
$ws.Range("L97").Value = 9000
$ws.Range("M97").Value = 8500
$ws.Range("P97").Value = 142

$ws.Range("D98").Value = 44357

$ws.Range("D99").Value = 44558
$ws.Range("J99").Value = 250
$ws.Range("K99").Value = 8000
$ws.Range("L99").Value = 8500
$ws.Range("M99").Value = 8300
$ws.Range("O99").Value = "Provincia de Limarí"
$ws.Range("P99").Value = 138

$ws.Range("D100").Value = 44552
$ws.Range("J100").Value = 100
$ws.Range("K100").Value = 7000
$ws.Range("L100").Value = 8000
$ws.Range("M100").Value = 7500
$ws.Range("N100").Value = "`$/caja 60 unidades"
$ws.Range("O100").Value = "Región de Arica y Parinacota"
$ws.Range("P100").Value = 125
$ws.Range("Q100").Value = 60

$ws.Range("D101").Value = 44488
$ws.Range("K101").Value = 7000
$ws.Range("L101").Value = 7500
$ws.Range("M101").Value = 7250
$ws.Range("P101").Value = 121

$ws.Range("D102").Value = 44166
$ws.Range("J102").Value = 100
$ws.Range("K102").Value = 6500
$ws.Range("L102").Value = 7000
$ws.Range("M102").Value = 6750
$ws.Range("N102").Value = "`$/caja 60 unidades"
$ws.Range("O102").Value = "Región de Arica y Parinacota"
$ws.Range("P102").Value = 112
$ws.Range("Q102").Value = 60

$ws.Range("D103").Value = 44316
$ws.Range("K103").Value = 9000
$ws.Range("L103").Value = 10000
$ws.Range("M103").Value = 9500
$ws.Range("P103").Value = 158

$ws.Range("D104").Value = 44568
$ws.Range("K104").Value = 7000
$ws.Range("L104").Value = 7500
$ws.Range("M104").Value = 7250
$ws.Range("P104").Value = 121

$ws.Range("D105").Value = 44656
$ws.Range("J105").Value = 240
$ws.Range("K105").Value = 12000
$ws.Range("L105").Value = 13000
$ws.Range("M105").Value = 12500
$ws.Range("P105").Value = 208

$ws.Range("D106").Value = 44586
$ws.Range("J106").Value = 90
$ws.Range("K106").Value = 11000
$ws.Range("L106").Value = 12000
$ws.Range("M106").Value = 11444
$ws.Range("P106").Value = 191

$ws.Range("D107").Value = 44469
$ws.Range("J107").Value = 100
$ws.Range("K107").Value = 16000
$ws.Range("L107").Value = 17000
$ws.Range("M107").Value = 16500
$ws.Range("O107").Value = "Región de Arica y Parinacota"
$ws.Range("P107").Value = 275

$ws.Range("D108").Value = 44463
$ws.Range("K108").Value = 15000
$ws.Range("L108").Value = 15500
$ws.Range("M108").Value = 15250
$ws.Range("P108").Value = 254

$ws.Range("D109").Value = 44243
$ws.Range("K109").Value = 11000
$ws.Range("L109").Value = 12000
$ws.Range("M109").Value = 11500
$ws.Range("P109").Value = 192

$ws.Range("D110").Value = 44519
$ws.Range("J110").Value = 450
$ws.Range("K110").Value = 6500
$ws.Range("L110").Value = 7000
$ws.Range("M110").Value = 6778
$ws.Range("O110").Value = "Región de Arica y Parinacota"
$ws.Range("P110").Value = 113

$ws.Range("D111").Value = 44505
$ws.Range("J111").Value = 300
$ws.Range("K111").Value = 6500
$ws.Range("L111").Value = 7000
$ws.Range("M111").Value = 6750
$ws.Range("N111").Value = "`$/caja 80 unidades"
$ws.Range("O111").Value = "Región del Maule"
$ws.Range("P111").Value = 84
$ws.Range("Q111").Value = 80

$ws.Range("D112").Value = 44372

$ws.Range("D113").Value = 44671
$ws.Range("J113").Value = 180
$ws.Range("K113").Value = 13000
$ws.Range("L113").Value = 14000
$ws.Range("M113").Value = 13444
$ws.Range("O113").Value = "Región Metropolitana"
$ws.Range("P113").Value = 224

$ws.Range("D114").Value = 44474
$ws.Range("K114").Value = 19000
$ws.Range("L114").Value = 20000
$ws.Range("M114").Value = 19500
$ws.Range("P114").Value = 325

$ws.Range("D115").Value = 44631
$ws.Range("J115").Value = 220
$ws.Range("K115").Value = 17000
$ws.Range("L115").Value = 18000
$ws.Range("M115").Value = 17545
$ws.Range("P115").Value = 292

$ws.Range("D116").Value = 44672
$ws.Range("K116").Value = 14000
$ws.Range("L116").Value = 15000
$ws.Range("M116").Value = 14455
$ws.Range("P116").Value = 241

$ws.Range("D117").Value = 44365
$ws.Range("J117").Value = 100
$ws.Range("K117").Value = 13000
$ws.Range("L117").Value = 14000
$ws.Range("M117").Value = 13500
$ws.Range("P117").Value = 225

$ws.Range("D118").Value = 44427
$ws.Range("J118").Value = 100
$ws.Range("K118").Value = 14000
$ws.Range("L118").Value = 15000
$ws.Range("M118").Value = 14500
$ws.Range("P118").Value = 242

$ws.Range("D119").Value = 44565
$ws.Range("K119").Value = 7000
$ws.Range("L119").Value = 8000
$ws.Range("M119").Value = 7500
$ws.Range("P119").Value = 125

$ws.Range("D120").Value = 44447
$ws.Range("K120").Value = 16000
$ws.Range("L120").Value = 17000
$ws.Range("M120").Value = 16500
$ws.Range("P120").Value = 275

$ws.Range("D121").Value = 44523
$ws.Range("K121").Value = 6500
$ws.Range("L121").Value = 7000
$ws.Range("M121").Value = 6750
$ws.Range("P121").Value = 112

$ws.Range("D122").Value = 44343
$ws.Range("J122").Value = 100
$ws.Range("L122").Value = 11000
$ws.Range("M122").Value = 10500
$ws.Range("O122").Value = "Región de Arica y Parinacota"
$ws.Range("P122").Value = 175

$ws.Range("D123").Value = 44462
$ws.Range("K123").Value = 14500
$ws.Range("L123").Value = 15000
$ws.Range("M123").Value = 14750
$ws.Range("P123").Value = 246

$ws.Range("D124").Value = 44377
$ws.Range("K124").Value = 11000
$ws.Range("L124").Value = 12000
$ws.Range("M124").Value = 11500
$ws.Range("P124").Value = 192

$ws.Range("D125").Value = 44664
$ws.Range("K125").Value = 19000
$ws.Range("L125").Value = 20000
$ws.Range("M125").Value = 19500
$ws.Range("P125").Value = 325

$ws.Range("D126").Value = 44600
$ws.Range("J126").Value = 250
$ws.Range("K126").Value = 7000
$ws.Range("L126").Value = 8000
$ws.Range("M126").Value = 7520
$ws.Range("O126").Value = "Región del Maule"
$ws.Range("P126").Value = 125

$ws.Range("D127").Value = 44323
$ws.Range("J127").Value = 100
$ws.Range("K127").Value = 9000
$ws.Range("L127").Value = 10000
$ws.Range("M127").Value = 9500
$ws.Range("O127").Value = "Región de Arica y Parinacota"
$ws.Range("P127").Value = 158

$ws.Range("D128").Value = 44602
$ws.Range("J128").Value = 220
$ws.Range("K128").Value = 10000
$ws.Range("L128").Value = 11000
$ws.Range("M128").Value = 10545
$ws.Range("N128").Value = "`$/caja 80 unidades"
$ws.Range("O128").Value = "Región del Maule"
$ws.Range("P128").Value = 132
$ws.Range("Q128").Value = 80
